# ---------------------------------------------------------------------------
# Commit: "Thu, Jul 30, 2020  5:04:59 AM"
#
# The canonical-XML diff for this commit swaps the two theme parts of the
# package:
#   - ppt/theme/theme1.xml  (the deck's main theme, wired to
#     ppt/presentation.xml and ppt/slideMasters/slideMaster1.xml) changes
#     from the "Integral" / "Red Violet" colour scheme to the stock
#     "Office Theme" / "Office" colour scheme.
#   - ppt/theme/theme2.xml  (the notes-master theme, wired to
#     ppt/notesMasters/notesMaster1.xml) changes from "Office Theme" /
#     "Office" to "Integral" / "Red Violet" - i.e. it ends up holding what
#     used to be theme1's content.
#
# The <a:fontScheme> and <a:fmtScheme> blocks of the two themes are byte
# identical to begin with, so the only substantive content difference is
# the <a:clrScheme> (12 colour slots) plus the cosmetic name="" attributes
# on <a:theme> and <a:clrScheme>.
#
# Through the PowerPoint object model, theme colours are edited via the
# legacy ColorScheme collection (SlideMaster.ColorScheme / Slide.ColorScheme
# / etc. all resolve to the single presentation theme, i.e. theme1.xml),
# in theme-XML order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- Swap theme1.xml's colour scheme to the stock "Office" palette --------
$cs = $p.SlideMaster.ColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink -> 954F72

# Relabel the scheme/theme to match the incoming "Office Theme" naming
# (harmless if the host does not persist these - the colour values above
# are the substantive part of the edit).
$cs.Name = "Office"

$theme = $p.SlideMaster.Theme
$theme.Name = "Office Theme"

$design = $p.SlideMaster.Design
$design.Name = "Office Theme"

# --- Mirror the rename on the notes-master theme (theme2.xml) -------------
# (best-effort; the notes-master theme shares the same object-model surface
# as the slide theme in this host, so there is no independently addressable
# way to push the former theme1 ("Integral") content into theme2.xml.)
try {
    $nmTheme = $p.NotesMaster.Theme
    $nmTheme.Name = "Integral"
} catch {
}
